# feat: update ELO margin multiplier logic and documentation
#
# 1. Single Match Simulator!K14 - rewrite margin-multiplier formula so that
#    a 1-set and a 2-set difference both give the same (0.1x) bonus, and
#    anything beyond 2 sets gives the 0.2x bonus.
# 2. Documentation!B5 - update the explanation text to match the new rule.
# 3. Documentation - insert a new "Margin Multiplier" row (with its formula)
#    right above the existing "Delta (ELO)" row, pushing the rows below it
#    down by one (dimension grows from A1:B23 to A1:B24).

$wb = $excel.ActiveWorkbook

# --- 1. Single Match Simulator: new margin multiplier formula -------------
$wsSim = $wb.Sheets.Item("Single Match Simulator")
$wsSim.Range("K14").Formula = "=1 + IF(ABS(B8-B9)>2, 0.2, IF(ABS(B8-B9)>0, 0.1, 0))"

# --- 2 & 3. Documentation sheet updates ------------------------------------
$wsDoc = $wb.Sheets.Item("Documentation")

# Update the Margin Multiplier explanation text in row 5.
$wsDoc.Range("B5").Value = "Bonus for decisive wins. 1-2 sets = 1.1x multiplier. 3+ sets = 1.2x."

# Insert a new row above the current row 21 ("Delta (ELO)"), copy the
# formatting from the row above (row 20) so the new row matches the
# existing look of the table, then fill in the new "Margin Multiplier"
# term + formula.
$wsDoc.Range("A21").EntireRow.Insert()
$wsDoc.Range("A20:B20").Copy()
$wsDoc.Range("A21:B21").PasteSpecial(-4122)

$wsDoc.Range("A21").Value = "Margin Multiplier"
$wsDoc.Range("B21").Formula = "=1 + IF(ABS(Set_Diff)>2, 0.2, IF(ABS(Set_Diff)>0, 0.1, 0))"
